# Applies the edits captured by the target diff:
#  - Column B values "a"/"b"/"c" -> "A"/"B"/"C" (rows 2-13)
#  - Active cell / selection moves from A14 to B13
#  - Sheet default column width nudges from 11.55078125 to 11.5703125
#  - Header/footer font style name changes from "Regular" to "Normal"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Uppercase the single-letter "Configuracion" column values (B2:B13).
$ws.Range("B2").Value = "A"
$ws.Range("B3").Value = "A"
$ws.Range("B4").Value = "A"
$ws.Range("B5").Value = "A"
$ws.Range("B6").Value = "B"
$ws.Range("B7").Value = "B"
$ws.Range("B8").Value = "B"
$ws.Range("B9").Value = "B"
$ws.Range("B10").Value = "C"
$ws.Range("B11").Value = "C"
$ws.Range("B12").Value = "C"
$ws.Range("B13").Value = "C"

# Move the active selection to B13.
$ws.Range("B13").Select()

# Widen the sheet's default (standard) column width slightly.
$ws.StandardWidth = 11.5703125

# Switch the header/footer font style name from "Regular" to "Normal".
$ws.PageSetup.CenterHeader = '&"Times New Roman,Normal"&12&A'
$ws.PageSetup.CenterFooter = '&"Times New Roman,Normal"&12Page &P'
